# Auto-generated cell updates applying the Marilith_Profits diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1893.75
$ws.Range("I2").Value = 1475
$ws.Range("J2").Value = 2312.5
$ws.Range("K2").Value = 1475
$ws.Range("L2").Value = 2312.5
$ws.Range("M2").Value = -1362
$ws.Range("N2").Value = -2538.5

$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()

$ws.Range("H38").Value = 35.5
$ws.Range("I38").Value = 35.5
$ws.Range("K38").Value = 106.5
$ws.Range("M38").Value = 265.5

$ws.Range("H40").Value = 4478.6313
$ws.Range("I40").Value = 3450
$ws.Range("J40").Value = 4599.647
$ws.Range("K40").Value = 3450
$ws.Range("L40").Value = 4599.647
$ws.Range("M40").Value = -3275
$ws.Range("N40").Value = -4949.647

$ws.Range("H103").Value = 1467
$ws.Range("I103").Value = 1467
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 4401
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -3815
$ws.Range("N103").ClearContents()

$ws.Range("H107").Value = 792.84375
$ws.Range("I107").Value = 714.5
$ws.Range("J107").Value = 1132.3334
$ws.Range("K107").Value = 714.5
$ws.Range("L107").Value = 1132.3334
$ws.Range("M107").Value = 1205.5
$ws.Range("N107").Value = -4972.3334

$ws.Range("H137").Value = 2336.8333
$ws.Range("I137").Value = 1782.4445
$ws.Range("K137").Value = 5347.333500000001
$ws.Range("M137").Value = -2797.333500000001

$ws.Range("H138").Value = 4194.3887
$ws.Range("J138").Value = 4194.3887
$ws.Range("L138").Value = 12583.1661
$ws.Range("N138").Value = -22863.1661

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 2000
$ws.Range("J4").Value = 2000
$ws.Range("L4").Value = 2000
$ws.Range("N4").Value = -2232

$ws.Range("H5").Value = 96
$ws.Range("I5").Value = 96
$ws.Range("K5").Value = 96
$ws.Range("M5").Value = 16

$ws.Range("H35").Value = 1775
$ws.Range("I35").Value = 1775
$ws.Range("K35").Value = 1775
$ws.Range("M35").Value = -1369

$ws.Range("H39").Value = 4250
$ws.Range("I39").Value = 4250
$ws.Range("K39").Value = 4250
$ws.Range("M39").Value = -3730

$ws.Range("H74").Value = 1847.2
$ws.Range("I74").Value = 1444.2778
$ws.Range("K74").Value = 1444.2778
$ws.Range("M74").Value = -570.2778000000001

$ws.Range("H77").Value = 1847.2
$ws.Range("I77").Value = 1444.2778
$ws.Range("K77").Value = 7221.389
$ws.Range("M77").Value = -2853.389

$ws.Range("H122").Value = 3495.25
$ws.Range("I122").Value = 3495.25
$ws.Range("K122").Value = 10485.75
$ws.Range("M122").Value = -8035.75

$ws.Range("H132").Value = 1241.3334
$ws.Range("I132").Value = 1241.3334
$ws.Range("K132").Value = 3724.0002
$ws.Range("M132").Value = -1194.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 96
$ws.Range("I4").Value = 96
$ws.Range("K4").Value = 96
$ws.Range("M4").Value = 19

$ws.Range("H86").Value = 2000
$ws.Range("I86").Value = 1000
$ws.Range("J86").Value = 3000
$ws.Range("K86").Value = 1000
$ws.Range("L86").Value = 3000
$ws.Range("M86").Value = 123
$ws.Range("N86").Value = -5246

$ws.Range("H89").Value = 2000
$ws.Range("I89").Value = 1000
$ws.Range("J89").Value = 3000
$ws.Range("K89").Value = 5000
$ws.Range("L89").Value = 15000
$ws.Range("M89").Value = 616
$ws.Range("N89").Value = -26232

$ws.Range("H134").Value = 12322.758
$ws.Range("I134").Value = 9110.134
$ws.Range("J134").Value = 14999.944
$ws.Range("K134").Value = 27330.402
$ws.Range("L134").Value = 44999.83199999999
$ws.Range("M134").Value = -24795.402
$ws.Range("N134").Value = -50069.83199999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws.Range("H122").Value = 2914.6316
$ws.Range("I122").Value = 2914.6316
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8743.8948
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -6293.8948
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 1484
$ws.Range("I132").Value = 1430.8
$ws.Range("K132").Value = 4292.4
$ws.Range("M132").Value = -1762.4

$ws.Range("H134").Value = 4651.125
$ws.Range("I134").Value = 4651.125
$ws.Range("K134").Value = 13953.375
$ws.Range("M134").Value = -11418.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 170.25
$ws.Range("I2").Value = 80.333336
$ws.Range("K2").Value = 482.000016
$ws.Range("M2").Value = -369.000016

$ws.Range("H29").Value = 383.33334
$ws.Range("J29").Value = 383.33334
$ws.Range("L29").Value = 1150.00002
$ws.Range("N29").Value = -1704.00002

$ws.Range("H37").Value = 67500
$ws.Range("J37").Value = 67500
$ws.Range("L37").Value = 202500
$ws.Range("N37").Value = -202724

$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()

$ws.Range("H80").Value = 502.5
$ws.Range("I80").Value = 2
$ws.Range("K80").Value = 6
$ws.Range("M80").Value = 930

$ws.Range("H83").Value = 502.5
$ws.Range("I83").Value = 2
$ws.Range("K83").Value = 18
$ws.Range("M83").Value = 4662

$ws.Range("H97").Value = 1286.7693
$ws.Range("I97").Value = 1682.75
$ws.Range("J97").Value = 1110.7778
$ws.Range("K97").Value = 5048.25
$ws.Range("L97").Value = 3332.3334
$ws.Range("M97").Value = -4552.25
$ws.Range("N97").Value = -4324.3334

$ws.Range("H113").Value = 2320
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 2320
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 6960
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -11300

$ws.Range("H128").Value = 342000
$ws.Range("I128").Value = 342000
$ws.Range("K128").Value = 1026000
$ws.Range("M128").Value = -1021020

$ws.Range("H131").Value = 13935.647
$ws.Range("I131").Value = 23468
$ws.Range("J131").Value = 3211.75
$ws.Range("K131").Value = 70404
$ws.Range("L131").Value = 9635.25
$ws.Range("M131").Value = -65364
$ws.Range("N131").Value = -19715.25

$ws.Range("H137").Value = 10000
$ws.Range("I137").Value = 1500
$ws.Range("J137").Value = 18500
$ws.Range("K137").Value = 4500
$ws.Range("L137").Value = 55500
$ws.Range("M137").Value = 600
$ws.Range("N137").Value = -65700

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2873
$ws.Range("I102").Value = 2873
$ws.Range("K102").Value = 2873
$ws.Range("M102").Value = -1251

$ws.Range("H113").Value = 1057.6666
$ws.Range("I113").Value = 1057.6666
$ws.Range("K113").Value = 1057.6666
$ws.Range("M113").Value = 1112.3334

$ws.Range("H126").Value = 1288.8
$ws.Range("I126").Value = 1311
$ws.Range("J126").Value = 1200
$ws.Range("K126").Value = 3933
$ws.Range("L126").Value = 3600
$ws.Range("M126").Value = -1463
$ws.Range("N126").Value = -8540

$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4388.5
$ws.Range("I122").Value = 4272
$ws.Range("K122").Value = 12816
$ws.Range("M122").Value = -10366

$ws.Range("H139").Value = 49998
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1400
$ws.Range("I96").Value = 1000
$ws.Range("J96").Value = 1800
$ws.Range("K96").Value = 1000
$ws.Range("L96").Value = 1800
$ws.Range("M96").Value = 373
$ws.Range("N96").Value = -4546
